# Auto-generated script to apply odds updates to sheet1 (row 2-7, columns F..AO)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 6).Value = 1.73   # F2
$ws.Cells.Item(2, 7).Value = 1.83   # G2
$ws.Cells.Item(2, 8).Value = 4.9   # H2
$ws.Cells.Item(2, 9).Value = 5.7   # I2
$ws.Cells.Item(2, 10).Value = 3.75   # J2
$ws.Cells.Item(2, 11).Value = 4.2   # K2
$ws.Cells.Item(2, 12).Value = 1.4   # L2
$ws.Cells.Item(2, 13).Value = 1.07   # M2
$ws.Cells.Item(2, 14).Value = 3.7   # N2
$ws.Cells.Item(2, 15).Value = 1.32   # O2
$ws.Cells.Item(2, 16).Value = 2   # P2
$ws.Cells.Item(2, 17).Value = 1.88   # Q2
$ws.Cells.Item(2, 18).Value = 1.37   # R2
$ws.Cells.Item(2, 19).Value = 3.25   # S2
$ws.Cells.Item(2, 21).Value = 2   # U2
$ws.Cells.Item(2, 22).Value = 1.21   # V2
$ws.Cells.Item(2, 23).Value = 2.2   # W2
$ws.Cells.Item(2, 24).Value = 15.5   # X2
$ws.Cells.Item(2, 25).Value = 19.5   # Y2
$ws.Cells.Item(2, 26).Value = 980   # Z2
$ws.Cells.Item(2, 28).Value = 8.800000000000001   # AB2
$ws.Cells.Item(2, 29).Value = 9.199999999999999   # AC2
$ws.Cells.Item(2, 30).Value = 23   # AD2
$ws.Cells.Item(2, 32).Value = 12   # AF2
$ws.Cells.Item(2, 33).Value = 10   # AG2
$ws.Cells.Item(2, 34).Value = 42   # AH2
$ws.Cells.Item(2, 36).Value = 21   # AJ2
$ws.Cells.Item(2, 37).Value = 38   # AK2
$ws.Cells.Item(2, 38).Value = 170   # AL2
$ws.Cells.Item(2, 41).Value = 600   # AO2

# Row 3
$ws.Cells.Item(3, 7).Value = 1.64   # G3
$ws.Cells.Item(3, 8).Value = 6   # H3
$ws.Cells.Item(3, 11).Value = 4.6   # K3
$ws.Cells.Item(3, 12).Value = 1.35   # L3
$ws.Cells.Item(3, 14).Value = 4.5   # N3
$ws.Cells.Item(3, 17).Value = 1.73   # Q3
$ws.Cells.Item(3, 18).Value = 1.47   # R3
$ws.Cells.Item(3, 19).Value = 2.92   # S3
$ws.Cells.Item(3, 20).Value = 1.79   # T3
$ws.Cells.Item(3, 21).Value = 2.1   # U3
$ws.Cells.Item(3, 22).Value = 1.17   # V3
$ws.Cells.Item(3, 23).Value = 2.56   # W3
$ws.Cells.Item(3, 25).Value = 26   # Y3
$ws.Cells.Item(3, 26).Value = 55   # Z3
$ws.Cells.Item(3, 27).Value = 180   # AA3
$ws.Cells.Item(3, 28).Value = 10.5   # AB3
$ws.Cells.Item(3, 29).Value = 10   # AC3
$ws.Cells.Item(3, 31).Value = 240   # AE3
$ws.Cells.Item(3, 33).Value = 9.6   # AG3
$ws.Cells.Item(3, 34).Value = 22   # AH3
$ws.Cells.Item(3, 36).Value = 15.5   # AJ3
$ws.Cells.Item(3, 40).Value = 7.8   # AN3
$ws.Cells.Item(3, 41).Value = 90   # AO3

# Row 4
$ws.Cells.Item(4, 7).Value = 2.06   # G4
$ws.Cells.Item(4, 8).Value = 3.65   # H4
$ws.Cells.Item(4, 9).Value = 4.1   # I4
$ws.Cells.Item(4, 10).Value = 3.85   # J4
$ws.Cells.Item(4, 12).Value = 1.34   # L4
$ws.Cells.Item(4, 14).Value = 4.4   # N4
$ws.Cells.Item(4, 15).Value = 1.23   # O4
$ws.Cells.Item(4, 16).Value = 2.2   # P4
$ws.Cells.Item(4, 17).Value = 1.71   # Q4
$ws.Cells.Item(4, 20).Value = 1.63   # T4
$ws.Cells.Item(4, 21).Value = 2.32   # U4
$ws.Cells.Item(4, 23).Value = 1.94   # W4
$ws.Cells.Item(4, 26).Value = 75   # Z4
$ws.Cells.Item(4, 27).Value = 1000   # AA4
$ws.Cells.Item(4, 31).Value = 190   # AE4
$ws.Cells.Item(4, 33).Value = 11   # AG4
$ws.Cells.Item(4, 34).Value = 17   # AH4
$ws.Cells.Item(4, 35).Value = 170   # AI4
$ws.Cells.Item(4, 36).Value = 55   # AJ4
$ws.Cells.Item(4, 37).Value = 21   # AK4
$ws.Cells.Item(4, 38).Value = 70   # AL4
$ws.Cells.Item(4, 39).Value = 330   # AM4
$ws.Cells.Item(4, 41).Value = 210   # AO4

# Row 5
$ws.Cells.Item(5, 6).Value = 2.3   # F5
$ws.Cells.Item(5, 10).Value = 3.05   # J5
$ws.Cells.Item(5, 11).Value = 3.45   # K5
$ws.Cells.Item(5, 12).Value = 1.53   # L5
$ws.Cells.Item(5, 14).Value = 2.78   # N5
$ws.Cells.Item(5, 15).Value = 1.46   # O5
$ws.Cells.Item(5, 16).Value = 1.6   # P5
$ws.Cells.Item(5, 17).Value = 2.38   # Q5
$ws.Cells.Item(5, 22).Value = 1.34   # V5
$ws.Cells.Item(5, 23).Value = 1.7   # W5
$ws.Cells.Item(5, 24).Value = 10.5   # X5
$ws.Cells.Item(5, 25).Value = 1000   # Y5
$ws.Cells.Item(5, 29).Value = 7.4   # AC5
$ws.Cells.Item(5, 30).Value = 16.5   # AD5
$ws.Cells.Item(5, 32).Value = 1000   # AF5
$ws.Cells.Item(5, 39).Value = 170   # AM5
$ws.Cells.Item(5, 41).Value = 80   # AO5

# Row 6
$ws.Cells.Item(6, 6).Value = 2.54   # F6
$ws.Cells.Item(6, 8).Value = 3.1   # H6
$ws.Cells.Item(6, 9).Value = 3.6   # I6
$ws.Cells.Item(6, 11).Value = 3.25   # K6
$ws.Cells.Item(6, 12).Value = 1.6   # L6
$ws.Cells.Item(6, 13).Value = 1.13   # M6
$ws.Cells.Item(6, 15).Value = 1.56   # O6
$ws.Cells.Item(6, 17).Value = 2.66   # Q6
$ws.Cells.Item(6, 19).Value = 5.3   # S6
$ws.Cells.Item(6, 20).Value = 2.1   # T6
$ws.Cells.Item(6, 21).Value = 1.73   # U6
$ws.Cells.Item(6, 22).Value = 1.39   # V6
$ws.Cells.Item(6, 28).Value = 19.5   # AB6

# Row 7
$ws.Cells.Item(7, 6).Value = 3.55   # F7
$ws.Cells.Item(7, 7).Value = 4.2   # G7
$ws.Cells.Item(7, 8).Value = 2.16   # H7
$ws.Cells.Item(7, 11).Value = 3.5   # K7
$ws.Cells.Item(7, 12).Value = 1.54   # L7
$ws.Cells.Item(7, 13).Value = 1.09   # M7
$ws.Cells.Item(7, 14).Value = 2.84   # N7
$ws.Cells.Item(7, 15).Value = 1.46   # O7
$ws.Cells.Item(7, 16).Value = 1.6   # P7
$ws.Cells.Item(7, 17).Value = 2.34   # Q7
$ws.Cells.Item(7, 18).Value = 1.22   # R7
$ws.Cells.Item(7, 19).Value = 4.5   # S7
$ws.Cells.Item(7, 20).Value = 1.96   # T7
$ws.Cells.Item(7, 21).Value = 1.83   # U7
$ws.Cells.Item(7, 22).Value = 1.71   # V7
$ws.Cells.Item(7, 23).Value = 1.32   # W7
$ws.Cells.Item(7, 29).Value = 1000   # AC7
